$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.612.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.450.02"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.75"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +9.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.451.15"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.77%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.040.34"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.01"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +7.64%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.454.61"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.764.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +8.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.40"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.98%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +4.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.568"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.588.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.08"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.79%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.33%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -12.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.26"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.83%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.98"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0788"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +12.05%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.37"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.597.44"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.85%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.98"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.26"
